# "Generate Report for Handback" - mark the a.md / b.md rows as handed back,
# fill in the per-language handback file/date columns, and widen the
# columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61502a264614e7763592cf91f36bf31cf6c824b1/e2e"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn / de-de) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "$baseUrl/a.md", "", "", "a.md")
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 14:48:13"

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "$baseUrl/a.md", "", "", "a.md")
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 14:48:13"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "$baseUrl/a.md", "", "", "a.md")
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 14:48:26"

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "$baseUrl/a.md", "", "", "a.md")
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 14:48:26"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667
